$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 2340
$ws.Range("I20").Value = 2340
$ws.Range("K20").Value = 2340
$ws.Range("M20").Value = -2110

$ws.Range("H35").Value = 2340
$ws.Range("I35").Value = 2340
$ws.Range("K35").Value = 2340
$ws.Range("M35").Value = -1961

$ws.Range("H51").Value = 5867.278
$ws.Range("I51").Value = 3660
$ws.Range("J51").Value = 6308.7334
$ws.Range("K51").Value = 3660
$ws.Range("L51").Value = 6308.7334
$ws.Range("M51").Value = -3176
$ws.Range("N51").Value = -7276.7334

$ws.Range("H58").Value = 1786.6666
$ws.Range("I58").Value = 180
$ws.Range("J58").Value = 5000
$ws.Range("K58").Value = 540
$ws.Range("L58").Value = 15000
$ws.Range("M58").Value = -390
$ws.Range("N58").Value = -15300

$ws.Range("H64").Value = 4590
$ws.Range("J64").Value = 4928.5713
$ws.Range("L64").Value = 4928.5713
$ws.Range("N64").Value = -5424.5713

$ws.Range("H67").Value = 4590
$ws.Range("J67").Value = 4928.5713
$ws.Range("L67").Value = 4928.5713
$ws.Range("N67").Value = -6644.5713

$ws.Range("H74").Value = 3929.0476
$ws.Range("I74").Value = 3925.5557
$ws.Range("J74").Value = 3950
$ws.Range("K74").Value = 3925.5557
$ws.Range("L74").Value = 3950
$ws.Range("M74").Value = -2989.5557
$ws.Range("N74").Value = -5822

$ws.Range("H76").Value = 3260.6
$ws.Range("I76").Value = 3284.3333
$ws.Range("J76").Value = 3225
$ws.Range("K76").Value = 3284.3333
$ws.Range("L76").Value = 3225
$ws.Range("M76").Value = -2969.3333
$ws.Range("N76").Value = -3855

$ws.Range("H77").Value = 3929.0476
$ws.Range("I77").Value = 3925.5557
$ws.Range("J77").Value = 3950
$ws.Range("K77").Value = 19627.7785
$ws.Range("L77").Value = 19750
$ws.Range("M77").Value = -14947.7785
$ws.Range("N77").Value = -29110

$ws.Range("H79").Value = 3260.6
$ws.Range("I79").Value = 3284.3333
$ws.Range("J79").Value = 3225
$ws.Range("K79").Value = 3284.3333
$ws.Range("L79").Value = 3225
$ws.Range("M79").Value = -2192.3333
$ws.Range("N79").Value = -5409

$ws.Range("H97").Value = 2176
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 2176
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 6528
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -7520

$ws.Range("H100").Value = 1480.1875
$ws.Range("I100").Value = 932.3333
$ws.Range("J100").Value = 3123.75
$ws.Range("K100").Value = 932.3333
$ws.Range("L100").Value = 3123.75
$ws.Range("M100").Value = -391.3333
$ws.Range("N100").Value = -4205.75

$ws.Range("H137").Value = 2327591.2
$ws.Range("I137").Value = 2942450.8
$ws.Range("J137").Value = 4788.778
$ws.Range("K137").Value = 8827352.399999999
$ws.Range("L137").Value = 14366.334
$ws.Range("M137").Value = -8824802.399999999
$ws.Range("N137").Value = -19466.334

$ws.Range("H138").Value = 4208361
$ws.Range("I138").Value = 287622.28
$ws.Range("J138").Value = 23812054
$ws.Range("K138").Value = 862866.8400000001
$ws.Range("L138").Value = 71436162
$ws.Range("M138").Value = -857726.8400000001
$ws.Range("N138").Value = -71446442

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2499.2856
$ws.Range("I2").Value = 2411
$ws.Range("K2").Value = 2411
$ws.Range("M2").Value = -2298

$ws.Range("H74").Value = 7413888
$ws.Range("I74").Value = 8656386
$ws.Range("K74").Value = 8656386
$ws.Range("M74").Value = -8655512

$ws.Range("H77").Value = 7413888
$ws.Range("I77").Value = 8656386
$ws.Range("K77").Value = 43281930
$ws.Range("M77").Value = -43277562

$ws.Range("H88").Value = 7073.4287
$ws.Range("J88").Value = 7073.4287
$ws.Range("L88").Value = 7073.4287
$ws.Range("N88").Value = -7885.4287

$ws.Range("H91").Value = 7073.4287
$ws.Range("J91").Value = 7073.4287
$ws.Range("L91").Value = 7073.4287
$ws.Range("N91").Value = -9881.4287

$ws.Range("H97").Value = 1737021.9
$ws.Range("I97").Value = 2233087
$ws.Range("K97").Value = 2233087
$ws.Range("M97").Value = -2232591

$ws.Range("H102").Value = 8930401
$ws.Range("I102").Value = 10990570
$ws.Range("K102").Value = 10990570
$ws.Range("M102").Value = -10988948

$ws.Range("H110").Value = 1319.6
$ws.Range("I110").Value = 933.3333
$ws.Range("J110").Value = 1899
$ws.Range("K110").Value = 933.3333
$ws.Range("L110").Value = 1899
$ws.Range("M110").Value = 1111.6667
$ws.Range("N110").Value = -5989

$ws.Range("H116").Value = 2499.2856
$ws.Range("I116").Value = 2411
$ws.Range("K116").Value = 2411
$ws.Range("M116").Value = -117

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2499.2856
$ws.Range("I3").Value = 2411
$ws.Range("K3").Value = 2411
$ws.Range("M3").Value = -2297

$ws.Range("H105").Value = 45457256
$ws.Range("I105").Value = 71431544
$ws.Range("J105").Value = 2250
$ws.Range("K105").Value = 71431544
$ws.Range("L105").Value = 2250
$ws.Range("M105").Value = -71429797
$ws.Range("N105").Value = -5744

$ws.Range("H134").Value = 4843.2
$ws.Range("I134").Value = 5506
$ws.Range("J134").Value = 4677.5
$ws.Range("K134").Value = 16518
$ws.Range("L134").Value = 14032.5
$ws.Range("M134").Value = -13983
$ws.Range("N134").Value = -19102.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2855.3235
$ws.Range("I31").Value = 1527.3462
$ws.Range("J31").Value = 7171.25
$ws.Range("K31").Value = 1527.3462
$ws.Range("L31").Value = 7171.25
$ws.Range("M31").Value = -1232.3462
$ws.Range("N31").Value = -7761.25

$ws.Range("H34").Value = 2855.3235
$ws.Range("I34").Value = 1527.3462
$ws.Range("J34").Value = 7171.25
$ws.Range("K34").Value = 1527.3462
$ws.Range("L34").Value = 7171.25
$ws.Range("M34").Value = -1325.3462
$ws.Range("N34").Value = -7575.25

$ws.Range("H36").Value = 4774
$ws.Range("I36").Value = 4774
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 4774
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -4386
$ws.Range("N36").ClearContents()

$ws.Range("H40").Value = 4774
$ws.Range("I40").Value = 4774
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 4774
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -4614
$ws.Range("N40").ClearContents()

$ws.Range("H94").Value = 4983.25
$ws.Range("I94").Value = 23000
$ws.Range("J94").Value = 1379.9
$ws.Range("K94").Value = 23000
$ws.Range("L94").Value = 1379.9
$ws.Range("M94").Value = -22549
$ws.Range("N94").Value = -2281.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H43").Value = 2802
$ws.Range("I43").Value = 2802
$ws.Range("K43").Value = 8406
$ws.Range("M43").Value = -8292

$ws.Range("H86").Value = 2118.6
$ws.Range("J86").Value = 2523.25
$ws.Range("L86").Value = 7569.75
$ws.Range("N86").Value = -9941.75

$ws.Range("H89").Value = 2118.6
$ws.Range("J89").Value = 2523.25
$ws.Range("L89").Value = 22709.25
$ws.Range("N89").Value = -34565.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 34461.707
$ws.Range("I70").Value = 50159.047
$ws.Range("J70").Value = 5683.25
$ws.Range("K70").Value = 50159.047
$ws.Range("L70").Value = 5683.25
$ws.Range("M70").Value = -49889.047
$ws.Range("N70").Value = -6223.25

$ws.Range("H73").Value = 34461.707
$ws.Range("I73").Value = 50159.047
$ws.Range("J73").Value = 5683.25
$ws.Range("K73").Value = 50159.047
$ws.Range("L73").Value = 5683.25
$ws.Range("M73").Value = -49223.047
$ws.Range("N73").Value = -7555.25

$ws.Range("H80").Value = 3270
$ws.Range("J80").Value = 3666.8
$ws.Range("L80").Value = 3666.8
$ws.Range("N80").Value = -5662.8

$ws.Range("H83").Value = 3270
$ws.Range("J83").Value = 3666.8
$ws.Range("L83").Value = 18334
$ws.Range("N83").Value = -28318

$ws.Range("H113").Value = 1579.5333
$ws.Range("I113").Value = 1125.7142
$ws.Range("J113").Value = 1976.625
$ws.Range("K113").Value = 1125.7142
$ws.Range("L113").Value = 1976.625
$ws.Range("M113").Value = 1044.2858
$ws.Range("N113").Value = -6316.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 867.65216
$ws.Range("I22").Value = 637.5
$ws.Range("J22").Value = 1118.7273
$ws.Range("K22").Value = 637.5
$ws.Range("L22").Value = 1118.7273
$ws.Range("M22").Value = -342.5
$ws.Range("N22").Value = -1708.7273

$ws.Range("H27").Value = 867.65216
$ws.Range("I27").Value = 637.5
$ws.Range("J27").Value = 1118.7273
$ws.Range("K27").Value = 637.5
$ws.Range("L27").Value = 1118.7273
$ws.Range("M27").Value = -530.5
$ws.Range("N27").Value = -1332.7273

$ws.Range("H55").Value = 131.375
$ws.Range("I55").Value = 99.8
$ws.Range("J55").Value = 184
$ws.Range("K55").Value = 99.8
$ws.Range("L55").Value = 184
$ws.Range("M55").Value = 73.2
$ws.Range("N55").Value = -530

$ws.Range("H68").Value = 1632.6538
$ws.Range("I68").Value = 1593.1364
$ws.Range("J68").Value = 1850
$ws.Range("K68").Value = 1593.1364
$ws.Range("L68").Value = 1850
$ws.Range("M68").Value = -844.1364000000001
$ws.Range("N68").Value = -3348

$ws.Range("H71").Value = 1632.6538
$ws.Range("I71").Value = 1593.1364
$ws.Range("J71").Value = 1850
$ws.Range("K71").Value = 7965.682000000001
$ws.Range("L71").Value = 9250
$ws.Range("M71").Value = -4221.682000000001
$ws.Range("N71").Value = -16738

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 111666.82
$ws.Range("I132").Value = 118154.3
$ws.Range("J132").Value = 101640.73
$ws.Range("K132").Value = 354462.9
$ws.Range("L132").Value = 304922.19
$ws.Range("M132").Value = -351932.9
$ws.Range("N132").Value = -309982.19
